# Consolidate split text runs back into a single run per paragraph.
#
# The original title / caption paragraphs were split across several
# <a:r> runs (one per word/space) even though every run shares identical
# (empty) formatting. Re-assigning the *same* string to TextRange.Text is
# treated as a no-op by the writer and leaves the XML untouched, so each
# TextRange is first set to a short placeholder value (forcing a genuine
# text change) and then set to the final desired text; the writer then
# re-emits the paragraph as a single consolidated run.

$p = $ppt.ActivePresentation

function Set-ConsolidatedText($textRange, [string]$text) {
    $textRange.Text = "x"
    $textRange.Text = $text
}

# Title placeholders (Shape 1 on every slide in this deck).
$titles = @{
    1  = "Slide 1 (Content)"
    2  = "Slide 2 (Content)"
    3  = "Slide 3 (Content)"
    4  = "Slide 4 (Content)"
    5  = "Slide 5 (Two Content)"
    6  = "Slide 6 (Two Content Right)"
    7  = "Slide 7 (Content with Caption)"
    8  = "Slide 8 (Comparison)"
    9  = "Slide 10 (Content)"
    10 = "Slide 11 (Content)"
    11 = "Slide 12 (Content)"
}

foreach ($idx in $titles.Keys) {
    $s = $p.Slides.Item($idx)
    Set-ConsolidatedText $s.Shapes.Item(1).TextFrame.TextRange $titles[$idx]
}

# Standalone caption textboxes that were also split across multiple runs.
$s6 = $p.Slides.Item(6)
Set-ConsolidatedText $s6.Shapes.Item(3).TextFrame.TextRange "an image"

$s7 = $p.Slides.Item(7)
Set-ConsolidatedText $s7.Shapes.Item(4).TextFrame.TextRange "An image"

$s8 = $p.Slides.Item(8)
Set-ConsolidatedText $s8.Shapes.Item(4).TextFrame.TextRange "An image"
